{"js": "// The document body contains a single table of 20 rows x 5 columns.\n// Each cell holds one arithmetic expression (e.g. \"66+25=91\"). This\n// edit replaces every cell's expression with a new one, in document\n// (row-major) order, while leaving all formatting untouched.\n\n// New expression for every cell, in row-major (top-to-bottom,\n// left-to-right) order -- 20 rows of 5 values each.\nconst newGrid = [\n  [\"43+45=88\", \"88-8=80\", \"95-33=62\", \"63-38=25\", \"78+11=89\"],\n  [\"95-62=33\", \"54-23=31\", \"19+40=59\", \"38+9=47\", \"11+21=32\"],\n  [\"77-41=36\", \"46+51=97\", \"67-58=9\", \"25+42=67\", \"44+34=78\"],\n  [\"55+6=61\", \"79-13=66\", \"71+16=87\", \"44+3=47\", \"95-8=87\"],\n  [\"40-37=3\", \"40-28=12\", \"73+2=75\", \"82-2=80\", \"70-56=14\"],\n  [\"12+63=75\", \"98-8=90\", \"41+36=77\", \"13+19=32\", \"10+40=50\"],\n  [\"42-9=33\", \"69+3=72\", \"58-42=16\", \"30+33=63\", \"12+48=60\"],\n  [\"55+14=69\", \"42+39=81\", \"41+44=85\", \"1+82=83\", \"45-20=25\"],\n  [\"70-55=15\", \"5+34=39\", \"30-2=28\", \"57-28=29\", \"71+26=97\"],\n  [\"18+8=26\", \"65-51=14\", \"96-51=45\", \"50-10=40\", \"87+7=94\"],\n  [\"28+47=75\", \"72+23=95\", \"82-66=16\", \"51-8=43\", \"2+62=64\"],\n  [\"31-10=21\", \"8+12=20\", \"84-8=76\", \"20-1=19\", \"93-12=81\"],\n  [\"90+2=92\", \"22-21=1\", \"11+81=92\", \"96-68=28\", \"39+30=69\"],\n  [\"68-50=18\", \"88-48=40\", \"58-48=10\", \"97-87=10\", \"96-89=7\"],\n  [\"96-56=40\", \"60-5=55\", \"38-4=34\", \"62+37=99\", \"78+6=84\"],\n  [\"40+25=65\", \"53-51=2\", \"57-16=41\", \"13+52=65\", \"1+44=45\"],\n  [\"71+5=76\", \"27-11=16\", \"52+6=58\", \"53-29=24\", \"26+51=77\"],\n  [\"51-48=3\", \"61-28=33\", \"35-1=34\", \"8+84=92\", \"99-55=44\"],\n  [\"68-2=66\", \"16+49=65\", \"73+5=78\", \"92-68=24\", \"59-49=10\"],\n  [\"95-91=4\", \"86-22=64\", \"17-8=9\", \"68+15=83\", \"1+8=9\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Update each cell's text in place. Using getCell/getRange so that\n// existing run formatting (font, size, etc.) on the single run inside\n// each cell paragraph is preserved exactly, rather than replacing the\n// whole table (which would reset formatting to a default run).\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < newGrid[r].length; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.getRange().insertText(newGrid[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document body contains a single table of 20 rows x 5 columns.\n# Each cell holds one arithmetic expression (e.g. \"66+25=91\"). This\n# edit replaces every cell's expression with a new one, in document\n# (row-major) order, while leaving all formatting untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New expression for every cell, in row-major (top-to-bottom,\n# left-to-right) order -- 20 rows of 5 values each.\n$newGrid = @(\n    @(\"43+45=88\", \"88-8=80\", \"95-33=62\", \"63-38=25\", \"78+11=89\"),\n    @(\"95-62=33\", \"54-23=31\", \"19+40=59\", \"38+9=47\", \"11+21=32\"),\n    @(\"77-41=36\", \"46+51=97\", \"67-58=9\", \"25+42=67\", \"44+34=78\"),\n    @(\"55+6=61\", \"79-13=66\", \"71+16=87\", \"44+3=47\", \"95-8=87\"),\n    @(\"40-37=3\", \"40-28=12\", \"73+2=75\", \"82-2=80\", \"70-56=14\"),\n    @(\"12+63=75\", \"98-8=90\", \"41+36=77\", \"13+19=32\", \"10+40=50\"),\n    @(\"42-9=33\", \"69+3=72\", \"58-42=16\", \"30+33=63\", \"12+48=60\"),\n    @(\"55+14=69\", \"42+39=81\", \"41+44=85\", \"1+82=83\", \"45-20=25\"),\n    @(\"70-55=15\", \"5+34=39\", \"30-2=28\", \"57-28=29\", \"71+26=97\"),\n    @(\"18+8=26\", \"65-51=14\", \"96-51=45\", \"50-10=40\", \"87+7=94\"),\n    @(\"28+47=75\", \"72+23=95\", \"82-66=16\", \"51-8=43\", \"2+62=64\"),\n    @(\"31-10=21\", \"8+12=20\", \"84-8=76\", \"20-1=19\", \"93-12=81\"),\n    @(\"90+2=92\", \"22-21=1\", \"11+81=92\", \"96-68=28\", \"39+30=69\"),\n    @(\"68-50=18\", \"88-48=40\", \"58-48=10\", \"97-87=10\", \"96-89=7\"),\n    @(\"96-56=40\", \"60-5=55\", \"38-4=34\", \"62+37=99\", \"78+6=84\"),\n    @(\"40+25=65\", \"53-51=2\", \"57-16=41\", \"13+52=65\", \"1+44=45\"),\n    @(\"71+5=76\", \"27-11=16\", \"52+6=58\", \"53-29=24\", \"26+51=77\"),\n    @(\"51-48=3\", \"61-28=33\", \"35-1=34\", \"8+84=92\", \"99-55=44\"),\n    @(\"68-2=66\", \"16+49=65\", \"73+5=78\", \"92-68=24\", \"59-49=10\"),\n    @(\"95-91=4\", \"86-22=64\", \"17-8=9\", \"68+15=83\", \"1+8=9\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newGrid[$r - 1][$c - 1]\n    }\n}\n"}
